# Add a new results row (2024-07-10_v1) below the existing one, carry the
# highlighted "current result" formatting (green fill, style index 2) down
# to the new row, and clear that highlight from the old row since it is no
# longer the most recent entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (incl. green "latest result" fill) from row 2 so the
# new row 3 picks up the same style before we strip it off row 2.
$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)  # xlPasteFormats

# Populate the new row with the updated model results.
$ws.Range("A3").Value = "2024-07-10_v1"
$ws.Range("B3").Value = 0.9155
$ws.Range("C3").Value = 0.9379999999999999
$ws.Range("D3").Value = 0.893
$ws.Range("E3").Value = 0.8318426675174775
$ws.Range("F3").Value = 0.9703179999999999

# The previous row is no longer the latest result, so reset it back to the
# default (unstyled) look.
$ws.Range("A2:F2").Style = "Normal"
